$wb = $excel.ActiveWorkbook

# --- survey sheet: trim trailing spaces from several label cells ---
$survey = $wb.Worksheets.Item("survey")

$survey.Range("C5").Value  = "Text"
$survey.Range("C9").Value  = "Integer"
$survey.Range("C10").Value = "Decimal"
$survey.Range("C14").Value = "Date"
$survey.Range("C15").Value = "Time"
$survey.Range("C16").Value = "Date and time"
$survey.Range("C25").Value = "Geopoint"

# --- survey sheet: bump the thin spacer row (row 8) height ---
$survey.Rows.Item(8).RowHeight = 14.15

# --- survey sheet: move the view/selection from A1 to C9 ---
$survey.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$survey.Range("C9").Select()

# --- settings sheet: move the selection from A1 to B1 ---
$settings = $wb.Worksheets.Item("settings")
$settings.Activate()
$settings.Range("B1").Select()

# restore "survey" as the active/visible tab (it was selected before editing)
$survey.Activate()
